$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new doctor entry as the next row after the existing list.
$ws.Range("A32").Value = "Dr. Nitigya"
$ws.Range("B32").Value = "Timepass"

# Reflect where the user ended up editing: scrolled down and left the
# new row's second cell selected.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("B32").Select()
